$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $style = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $style
}

Set-TextValue $ws.Range("D2") "26.935.34"
$ws.Range("E2").Value = "  +0.96%  "
Set-TextValue $ws.Range("D3") "1.845.02"
$ws.Range("E3").Value = "  +0.92%  "
Set-TextValue $ws.Range("D4") "1.010"
$ws.Range("E4").Value = "  +0.27%  "
Set-TextValue $ws.Range("D5") "309.33"
$ws.Range("E5").Value = "  +0.34%  "
Set-TextValue $ws.Range("D6") "1.009"
$ws.Range("E6").Value = "  +0.19%  "
Set-TextValue $ws.Range("D7") "0.4781"
$ws.Range("E7").Value = "  +3.08%  "
$ws.Range("E8").Value = "  +1.73%  "
Set-TextValue $ws.Range("D9") "0.07209"
$ws.Range("E9").Value = "  +0.97%  "
Set-TextValue $ws.Range("D10") "0.9267"
$ws.Range("E10").Value = "  +2.59%  "
Set-TextValue $ws.Range("D11") "19.68"
$ws.Range("E11").Value = "  +1.30%  "
Set-TextValue $ws.Range("D12") "0.07686"
Set-TextValue $ws.Range("D13") "1.892.65"
$ws.Range("E13").Value = "  +3.13%  "
Set-TextValue $ws.Range("D14") "5.317"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("E15").Value = "  +0.77%  "
Set-TextValue $ws.Range("D16") "88.85"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("E17").Value = "  +0.24%  "
Set-TextValue $ws.Range("D18") "0.000008635"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("E19").Value = "  +0.13%  "
Set-TextValue $ws.Range("D20") "26.965.38"
$ws.Range("E20").Value = "  +0.92%  "
Set-TextValue $ws.Range("D21") "14.54"
$ws.Range("E21").Value = "  +2.45%  "
Set-TextValue $ws.Range("D22") "5.056"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("E24").Value = "  -0.36%  "
Set-TextValue $ws.Range("D25") "152.49"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  +1.23%  "
Set-TextValue $ws.Range("D27") "1.996"
$ws.Range("E27").Value = "  +0.85%  "
Set-TextValue $ws.Range("D28") "114.18"
$ws.Range("E28").Value = "  +0.21%  "
Set-TextValue $ws.Range("D29") "4.922"
Set-TextValue $ws.Range("D31") "3.316"
$ws.Range("E31").Value = "  +5.47%  "
Set-TextValue $ws.Range("D32") "1.172"
$ws.Range("E32").Value = "  +2.53%  "
Set-TextValue $ws.Range("D33") "0.7449"
$ws.Range("E33").Value = "  +1.40%  "
Set-TextValue $ws.Range("D34") "4.486"
$ws.Range("E34").Value = "  +0.93%  "
Set-TextValue $ws.Range("D35") "2.717"
$ws.Range("E35").Value = "  -0.12%  "
Set-TextValue $ws.Range("D36") "1.119"
$ws.Range("E36").Value = "  +4.06%  "
Set-TextValue $ws.Range("D37") "0.01957"
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("E38").Value = "  +2.42%  "
Set-TextValue $ws.Range("D39") "2.977"
Set-TextValue $ws.Range("D40") "0.5189"
$ws.Range("E40").Value = "  +2.23%  "
Set-TextValue $ws.Range("D41") "6.963"
$ws.Range("E41").Value = "  +1.23%  "
Set-TextValue $ws.Range("D42") "0.1508"
$ws.Range("E42").Value = "  +0.75%  "
Set-TextValue $ws.Range("D43") "8.190"
$ws.Range("E43").Value = "  +2.15%  "
Set-TextValue $ws.Range("D44") "10.57"
$ws.Range("E44").Value = "  +6.08%  "
Set-TextValue $ws.Range("D45") "0.4718"
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("E46").Value = "  +0.26%  "
Set-TextValue $ws.Range("D47") "101.40"
$ws.Range("E47").Value = "  +2.81%  "
Set-TextValue $ws.Range("D48") "1.600"
$ws.Range("E48").Value = "  +2.36%  "
Set-TextValue $ws.Range("D49") "65.88"
$ws.Range("E49").Value = "  +3.01%  "
Set-TextValue $ws.Range("D50") "0.06014"
$ws.Range("E50").Value = "  -0.17%  "
Set-TextValue $ws.Range("D51") "0.8849"
$ws.Range("E51").Value = "  +3.77%  "
